# Update column C ("Förändrad") for all data rows: 46074 -> 46075 (i.e. +1 day)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 149

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
